$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4243.7334
$ws.Range("I64").Value = 4325
$ws.Range("K64").Value = 4325
$ws.Range("M64").Value = -4077
$ws.Range("H67").Value = 4243.7334
$ws.Range("I67").Value = 4325
$ws.Range("K67").Value = 4325
$ws.Range("M67").Value = -3467
$ws.Range("H113").Value = 5519.4443
$ws.Range("I113").Value = 8444
$ws.Range("J113").Value = 3430.476
$ws.Range("K113").Value = 8444
$ws.Range("L113").Value = 3430.476
$ws.Range("M113").Value = -5190
$ws.Range("N113").Value = -9938.476000000001
$ws.Range("H116").Value = 189016.73
$ws.Range("I116").Value = 2517.5454
$ws.Range("J116").Value = 317234.94
$ws.Range("K116").Value = 2517.5454
$ws.Range("L116").Value = 317234.94
$ws.Range("M116").Value = 924.4546
$ws.Range("N116").Value = -324118.94
$ws.Range("H127").Value = 2198.9
$ws.Range("I127").Value = 1207.1111
$ws.Range("K127").Value = 3621.3333
$ws.Range("M127").Value = 1338.6667
$ws.Range("H132").Value = 37720344
$ws.Range("I132").Value = 46237570
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 138712710
$ws.Range("L132").Value = 3600
$ws.Range("M132").Value = -138710180
$ws.Range("N132").Value = -8660

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4731.881
$ws.Range("I32").Value = 4725.3413
$ws.Range("K32").Value = 4725.3413
$ws.Range("M32").Value = -4438.3413
$ws.Range("H45").Value = 944.8570999999999
$ws.Range("I45").Value = 906
$ws.Range("J45").Value = 996.6667
$ws.Range("K45").Value = 906
$ws.Range("L45").Value = 996.6667
$ws.Range("M45").Value = -529
$ws.Range("N45").Value = -1750.6667
$ws.Range("H61").Value = 5476.0835
$ws.Range("I61").Value = 6540.3887
$ws.Range("K61").Value = 6540.3887
$ws.Range("M61").Value = -6328.3887
$ws.Range("H102").Value = 166668140
$ws.Range("I102").Value = 166668140
$ws.Range("K102").Value = 166668140
$ws.Range("M102").Value = -166666518
$ws.Range("H132").Value = 3789681
$ws.Range("I132").Value = 4808979.5
$ws.Range("K132").Value = 14426938.5
$ws.Range("M132").Value = -14424408.5
$ws.Range("H136").Value = 5476.0835
$ws.Range("I136").Value = 6540.3887
$ws.Range("K136").Value = 19621.1661
$ws.Range("M136").Value = -17071.1661

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 9274107
$ws.Range("I134").Value = 10117041
$ws.Range("J134").Value = 1833.3334
$ws.Range("K134").Value = 30351123
$ws.Range("L134").Value = 5500.0002
$ws.Range("M134").Value = -30348588
$ws.Range("N134").Value = -10570.0002

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 6700
$ws.Range("I15").Value = 6500
$ws.Range("J15").Value = 6800
$ws.Range("K15").Value = 6500
$ws.Range("L15").Value = 6800
$ws.Range("M15").Value = -6330
$ws.Range("N15").Value = -7140
$ws.Range("H99").Value = 3287.4285
$ws.Range("I99").Value = 4102.4
$ws.Range("J99").Value = 1250
$ws.Range("K99").Value = 4102.4
$ws.Range("L99").Value = 1250
$ws.Range("M99").Value = -2604.4
$ws.Range("N99").Value = -4246
$ws.Range("H126").Value = 3287.4285
$ws.Range("I126").Value = 4102.4
$ws.Range("J126").Value = 1250
$ws.Range("K126").Value = 12307.2
$ws.Range("L126").Value = 3750
$ws.Range("M126").Value = -9837.199999999999
$ws.Range("N126").Value = -8690
$ws.Range("H132").Value = 8776377
$ws.Range("I132").Value = 9525586
$ws.Range("J132").Value = 35604.332
$ws.Range("K132").Value = 28576758
$ws.Range("L132").Value = 106812.996
$ws.Range("M132").Value = -28574228
$ws.Range("N132").Value = -111872.996
$ws.Range("H134").Value = 16448783
$ws.Range("I134").Value = 16668112
$ws.Range("J134").Value = 15626302
$ws.Range("K134").Value = 50004336
$ws.Range("L134").Value = 46878906
$ws.Range("M134").Value = -50001801
$ws.Range("N134").Value = -46883976

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 69373.125
$ws.Range("I104").Value = 1308.6666
$ws.Range("J104").Value = 85080.30499999999
$ws.Range("K104").Value = 3925.9998
$ws.Range("L104").Value = 255240.915
$ws.Range("M104").Value = -1304.9998
$ws.Range("N104").Value = -260482.915
$ws.Range("H113").Value = 2705918
$ws.Range("I113").Value = 720.25
$ws.Range("J113").Value = 3033820.8
$ws.Range("K113").Value = 2160.75
$ws.Range("L113").Value = 9101462.399999999
$ws.Range("M113").Value = 9.25
$ws.Range("N113").Value = -9105802.399999999
$ws.Range("H134").Value = 3745.1292
$ws.Range("I134").Value = 3811.7058
$ws.Range("J134").Value = 3664.2856
$ws.Range("K134").Value = 11435.1174
$ws.Range("L134").Value = 10992.8568
$ws.Range("M134").Value = -6365.117400000001
$ws.Range("N134").Value = -21132.8568

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 32990.8
$ws.Range("I70").Value = 47789.957
$ws.Range("K70").Value = 47789.957
$ws.Range("M70").Value = -47519.957
$ws.Range("H73").Value = 32990.8
$ws.Range("I73").Value = 47789.957
$ws.Range("K73").Value = 47789.957
$ws.Range("M73").Value = -46853.957
$ws.Range("H97").Value = 95238824
$ws.Range("I97").Value = 76923840
$ws.Range("K97").Value = 76923840
$ws.Range("M97").Value = -76923344
$ws.Range("H132").Value = 57146236
$ws.Range("I132").Value = 66667772
$ws.Range("J132").Value = 17030
$ws.Range("K132").Value = 200003316
$ws.Range("L132").Value = 51090
$ws.Range("M132").Value = -200000786
$ws.Range("N132").Value = -56150
$ws.Range("H17").Value = 400
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 400
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 400
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -736

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2270.5
$ws.Range("I7").Value = 1561
$ws.Range("J7").Value = 2980
$ws.Range("K7").Value = 1561
$ws.Range("L7").Value = 2980
$ws.Range("M7").Value = -1449
$ws.Range("N7").Value = -3204
$ws.Range("H11").Value = 70007
$ws.Range("J11").Value = 70007
$ws.Range("L11").Value = 70007
$ws.Range("N11").Value = -70287
$ws.Range("H100").Value = 2058.8333
$ws.Range("I100").Value = 1720.6
$ws.Range("J100").Value = 3750
$ws.Range("K100").Value = 1720.6
$ws.Range("L100").Value = 3750
$ws.Range("M100").Value = -1179.6
$ws.Range("N100").Value = -4832
$ws.Range("H122").Value = 50004640
$ws.Range("J122").Value = 50004640
$ws.Range("L122").Value = 150013920
$ws.Range("N122").Value = -150018820
$ws.Range("H126").Value = 2270.5
$ws.Range("I126").Value = 1561
$ws.Range("J126").Value = 2980
$ws.Range("K126").Value = 4683
$ws.Range("L126").Value = 8940
$ws.Range("M126").Value = -2213
$ws.Range("N126").Value = -13880
$ws.Range("H132").Value = 4547012
$ws.Range("I132").Value = 6250683.5
$ws.Range("K132").Value = 18752050.5
$ws.Range("M132").Value = -18749520.5
